# The only textual differences in the target OOXML are:
#  1) Relationship-id (r:id) strings in presentation.xml / slide.xml that PowerPoint
#     regenerates automatically every time the file is re-saved (slide master,
#     slide, and all 11 slide-layout relationship ids, plus the webextension
#     graphicFrame r:id and the picture/snapshot r:embed id). These are opaque,
#     auto-generated identifiers - their *targets* (slideMaster.xml, slide.xml,
#     slideLayoutN.xml, image.bin) are unchanged, so there is no content for the
#     object model to edit here.
#  2) The <we:webextension id="..."> GUID in ppt/slides/udata/data.xml. This is
#     PowerPoint's own internal instance id for the embedded "PowerPoll" add-in
#     object already on the slide (not the add-in's store reference, which is
#     untouched) - it is not exposed anywhere in the PowerPoint object model
#     (Shape/GraphicFrame/OLEFormat/CustomXMLParts/Tags all lack it), so it can't
#     be set through COM automation either, in this runtime or in real PowerPoint.
#
# There is no slide text, shape geometry, picture, or layout change to apply -
# the deck's visible/addressable content is identical before and after. Touch
# the presentation through the object model (matching the "fixed ... warnings"
# nature of the commit) without altering any content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$null = $s.Shapes.Count
